$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns retain exact text formatting (e.g. trailing zeros,
# thousand-dot separators) instead of being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "65.920.75"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").Value = "3.509.24"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "578.28"
$ws.Range("E5").Value = "  +4.65%  "
$ws.Range("D6").Value = "178.22"
$ws.Range("E6").Value = "  -6.25%  "
$ws.Range("E7").Value = "  +4.53%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("E10").Value = "  +3.98%  "
$ws.Range("D11").Value = "55.39"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("D13").Value = "9.23"
$ws.Range("E13").Value = "  -2.03%  "
$ws.Range("D14").Value = "4.075.50"
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("D15").Value = "3.508.63"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "18.36"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").Value = "12.05"
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("D19").Value = "65.882.41"
$ws.Range("E19").Value = "  -1.44%  "
$ws.Range("D20").Value = "1.00"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("D21").Value = "413.57"
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("E22").Value = "  +8.43%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "85.84"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "4.28"
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("D25").Value = "13.08"
$ws.Range("E25").Value = "  +9.69%  "
$ws.Range("D26").Value = "10.98"
$ws.Range("E26").Value = "  -1.34%  "
$ws.Range("E27").Value = "  -2.74%  "
$ws.Range("D28").Value = "9.07"
$ws.Range("E28").Value = "  +2.16%  "
$ws.Range("D29").Value = "30.36"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "624.10"
$ws.Range("E30").Value = "  -4.71%  "
$ws.Range("D31").Value = "6.47"
$ws.Range("E31").Value = "  -3.97%  "
$ws.Range("D32").Value = "11.64"
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").Value = "0.157"
$ws.Range("E34").Value = "  +14.40%  "
$ws.Range("D35").Value = "59.58"
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").Value = "0.0₃0797"
$ws.Range("E37").Value = "  -1.85%  "
$ws.Range("D38").Value = "37.29"
$ws.Range("E38").Value = "  -4.17%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "3.55"
$ws.Range("E39").Value = "  +5.94%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "3.302.02"
$ws.Range("E40").Value = "  +9.98%  "
$ws.Range("E41").Value = "  -3.53%  "
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").Value = "2.92"
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("E45").Value = "  -5.37%  "
$ws.Range("E46").Value = "  -4.37%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "0.133"
$ws.Range("E48").Value = "  +1.61%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "140.19"
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "8.57"
$ws.Range("E50").Value = "  -4.20%  "
$ws.Range("D51").Value = "2.31"
$ws.Range("E51").Value = "  -4.82%  "
